$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '96.796.44'
$ws.Range('E2').Value = '  -0.18%  '

# Row 3
$ws.Range('D3').Value = '3.713.32'
$ws.Range('E3').Value = '  +0.11%  '

# Row 4
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$ws.Range('D5').Value = '''236.28'
$ws.Range('E5').Value = '  -3.36%  '

# Row 6
$ws.Range('D6').Value = '''1.89'
$ws.Range('E6').Value = '  -0.88%  '

# Row 7
$ws.Range('D7').Value = '''651.45'
$ws.Range('E7').Value = '  -3.17%  '

# Row 8
$ws.Range('D8').Value = '''0.429'
$ws.Range('E8').Value = '  -0.46%  '

# Row 9
$ws.Range('E9').Value = '  -0.01%  '

# Row 10
$ws.Range('E10').Value = '  -6.26%  '

# Row 11
$ws.Range('D11').Value = '3.709.94'
$ws.Range('E11').Value = '  +0.14%  '

# Row 12
$ws.Range('B12').Value = 'ShibaInu'
$ws.Range('C12').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D12').Value = '''0.0000306'
$ws.Range('E12').Value = '  +13.80%  '

# Row 13
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').Value = '''44.14'
$ws.Range('E13').Value = '  -3.04%  '

# Row 14
$ws.Range('E14').Value = '  +0.20%  '

# Row 15
$ws.Range('E15').Value = '  +1.73%  '

# Row 16
$ws.Range('D16').Value = '4.403.21'
$ws.Range('E16').Value = '  +0.22%  '

# Row 17
$ws.Range('D17').Value = '96.502.27'

# Row 18
$ws.Range('D18').Value = '''8.84'
$ws.Range('E18').Value = '  -2.44%  '

# Row 19
$ws.Range('D19').Value = '3.723.45'
$ws.Range('E19').Value = '  +0.34%  '

# Row 20
$ws.Range('D20').Value = '''13.07'
$ws.Range('E20').Value = '  +0.44%  '

# Row 21
$ws.Range('D21').Value = '''18.69'
$ws.Range('E21').Value = '  +0.27%  '

# Row 22
$ws.Range('D22').Value = '''0.506'
$ws.Range('E22').Value = '  -6.53%  '

# Row 23
$ws.Range('D23').Value = '''521.14'
$ws.Range('E23').Value = '  +0.59%  '

# Row 24
$ws.Range('D24').Value = '''3.41'
$ws.Range('E24').Value = '  -1.73%  '

# Row 25
$ws.Range('E25').Value = '  +0.89%  '

# Row 26
$ws.Range('D26').Value = '''6.93'
$ws.Range('E26').Value = '  -0.98%  '

# Row 27
$ws.Range('D27').Value = '''101.56'
$ws.Range('E27').Value = '  -0.67%  '

# Row 28
$ws.Range('B28').Value = 'Hedera'
$ws.Range('C28').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D28').Value = '''0.187'
$ws.Range('E28').Value = '  +11.21%  '

# Row 29
$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').Value = '''13.34'
$ws.Range('E29').Value = '  +1.22%  '

# Row 30
$ws.Range('D30').Value = '''3.01'
$ws.Range('E30').Value = '  -3.13%  '

# Row 31
$ws.Range('D31').Value = '''12.14'
$ws.Range('E31').Value = '  -0.32%  '

# Row 32
$ws.Range('E32').Value = '  +0.01%  '

# Row 33
$ws.Range('B33').Value = 'Cronos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D33').Value = '''0.189'
$ws.Range('E33').Value = '  +0.97%  '

# Row 34
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').Value = '''1.87'
$ws.Range('E34').Value = '  +7.65%  '

# Row 35
$ws.Range('D35').Value = '''0.997'
$ws.Range('E35').Value = '  -0.27%  '

# Row 36
$ws.Range('D36').Value = '''32.32'
$ws.Range('E36').Value = '  -3.12%  '

# Row 37
$ws.Range('D37').Value = '''651.02'
$ws.Range('E37').Value = '  +5.69%  '

# Row 38
$ws.Range('D38').Value = '''0.590'
$ws.Range('E38').Value = '  -1.18%  '

# Row 39
$ws.Range('D39').Value = '''8.82'
$ws.Range('E39').Value = '  -0.43%  '

# Row 41
$ws.Range('E41').Value = '  -5.05%  '

# Row 42
$ws.Range('D42').Value = '''6.81'
$ws.Range('E42').Value = '  +9.45%  '

# Row 43
$ws.Range('E43').Value = '  +2.58%  '

# Row 44
$ws.Range('E44').Value = '  -2.83%  '

# Row 45
$ws.Range('E45').Value = '  -0.66%  '

# Row 46
$ws.Range('D46').Value = '''0.0452'
$ws.Range('E46').Value = '  -0.74%  '

# Row 47
$ws.Range('D47').Value = '''0.437'
$ws.Range('E47').Value = '  +0.72%  '

# Row 48
$ws.Range('E48').Value = '  -1.38%  '

# Row 49
$ws.Range('D49').Value = '''23.61'
$ws.Range('E49').Value = '  -0.03%  '

# Row 50
$ws.Range('D50').Value = '''8.50'
$ws.Range('E50').Value = '  -1.87%  '

# Row 51
$ws.Range('E51').Value = '  +0.98%  '
